$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: B2/C2 switch from the "-" placeholder text to numeric 0 ---
# (this also drops the now-unused "-" shared string entry on save)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# --- New header cells for the Area / Atotal columns and the summary block ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- New "Area" (G) column: per-segment cross-sectional area ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Atotal (H2) and the small summary block (J2:K2) ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Match the author's final selection / active cell ---
$ws.Range("J2:K2").Select() | Out-Null
